# Add two new benchmark sheets ("00-22 synchronous" and "00-30 multiprocessed")
# and wire them up into the "Main" summary sheet (rows 8 and 9), following the
# same layout used by the other "dim/count/time" worksheets already present
# in the workbook (e.g. "22-09 synchronous", "22-33 multiprocessed").

$wb = $excel.ActiveWorkbook

function Add-BenchmarkSheet {
    param($Name, $Dims, $Counts, $TimesSec, $TimesMin)

    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add($null, $lastSheet)
    $ws.Name = $Name

    $ws.Range("A1").Value = "dim"
    $ws.Range("B1").Value = "count"
    $ws.Range("C1").Value = "time, sec"
    $ws.Range("D1").Value = "time, min"
    $ws.Range("E1").Value = "avg matrix/min"
    $ws.Range("G1").Value = "total time spent, min"

    for ($i = 0; $i -lt $Dims.Length; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 1).Value = $Dims[$i]
        $ws.Cells.Item($row, 2).Value = $Counts[$i]
        $ws.Cells.Item($row, 3).Value = $TimesSec[$i]
        $ws.Cells.Item($row, 4).Value = $TimesMin[$i]
        $ws.Cells.Item($row, 5).Formula = "=60*B$row/C$row"
    }
    $ws.Range("G2").Formula = "=SUM(C2:C6)/60"
}

# --- "00-22 synchronous" -------------------------------------------------
Add-BenchmarkSheet "00-22 synchronous" `
    @(25, 50, 100, 200, 300) `
    @(1, 1, 1, 1, 1) `
    @(0.1970628000271972, 2.068523800000548, 29.42184979998274, 455.9260137999954, 2334.644614599994) `
    @(0, 0.03, 0.49, 7.6, 38.91) | Out-Null

# --- "00-30 multiprocessed" ----------------------------------------------
Add-BenchmarkSheet "00-30 multiprocessed" `
    @(25, 50, 100, 200, 300) `
    @(1, 1, 1, 1, 1) `
    @(1.78727249999065, 2.066149300022516, 5.765522300003795, 64.40812609999557, 390.2966042000044) `
    @(0.03, 0.03, 0.1, 1.07, 6.5) | Out-Null

# --- Wire the new sheets into the "Main" summary table --------------------
$main = $wb.Worksheets.Item("Main")

$main.Range("A8").Value = "00-22 synchronous"
$main.Range("B8").Formula = "='00-22 synchronous'!C2"
$main.Range("C8").Formula = "='00-22 synchronous'!C3"
$main.Range("D8").Formula = "='00-22 synchronous'!C4"
$main.Range("E8").Formula = "='00-22 synchronous'!C5"
$main.Range("F8").Formula = "='00-22 synchronous'!C6"

$main.Range("A9").Value = "00-30 multiprocessed"
$main.Range("B9").Formula = "='00-30 multiprocessed'!C2"
$main.Range("C9").Formula = "='00-30 multiprocessed'!C3"
$main.Range("D9").Formula = "='00-30 multiprocessed'!C4"
$main.Range("E9").Formula = "='00-30 multiprocessed'!C5"
$main.Range("F9").Formula = "='00-30 multiprocessed'!C6"
